$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-12-24 Tuesday"; new = "2024-12-25 Wednesday"},
    @{old = "736÷7="; new = "176÷7="},
    @{old = "538÷9="; new = "587÷9="},
    @{old = "350÷2="; new = "766÷5="},
    @{old = "886÷5="; new = "428÷5="},
    @{old = "319÷7="; new = "987÷8="},
    @{old = "519÷4="; new = "640÷8="},
    @{old = "445÷5="; new = "935÷4="},
    @{old = "163÷5="; new = "466÷9="},
    @{old = "607÷5="; new = "239÷7="},
    @{old = "951÷3="; new = "585÷6="},
    @{old = "976÷9="; new = "556÷4="},
    @{old = "382÷6="; new = "807÷2="},
    @{old = "722÷7="; new = "462÷3="},
    @{old = "656÷8="; new = "962÷6="},
    @{old = "261÷4="; new = "504÷6="},
    @{old = "884÷3="; new = "455÷8="},
    @{old = "865÷4="; new = "526÷4="},
    @{old = "346÷7="; new = "868÷8="},
    @{old = "997÷2="; new = "788÷2="},
    @{old = "841÷5="; new = "453÷5="},
    @{old = "289÷7="; new = "874÷3="},
    @{old = "314÷7="; new = "983÷7="},
    @{old = "274÷3="; new = "603÷6="},
    @{old = "515÷4="; new = "262÷7="},
    @{old = "963÷2="; new = "600÷8="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
